$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# The paragraph "Всего в электронную таблицу был занесены данные по 1000
# продуктов." is immediately followed by two more paragraphs: one that just
# says "Ссылка на таблицу:" and one that holds a centered hyperlink to
# "task14.xls". The edit removes that link entirely:
#   1. the "Ссылка на таблицу:" text collapses to a single trailing space
#      and its own paragraph mark is removed, so it becomes a second run
#      tacked onto the "...продуктов." paragraph;
#   2. the whole hyperlink paragraph (pilcrow included) is deleted outright.
# --------------------------------------------------------------------------

# 1) Turn the "Ссылка на таблицу:" run into a single space, in place, so it
#    keeps its own run formatting (font/size/color/lang) untouched.
$d.Content.Find.Execute("Ссылка на таблицу:", $false, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# 2) Merge that (now one-space) paragraph into the preceding
#    "...продуктов." paragraph by deleting the paragraph mark between them.
#    Paragraph.Range.End already points one past the paragraph's own pilcrow,
#    so the mark itself is the single character [End-1, End).
$targetRange = $d.Content
$targetRange.Find.Execute("Всего в электронную таблицу был занесены данные по 1000 продуктов.") | Out-Null
$totalsParagraph = $targetRange.Paragraphs(1)
$pilcrow = $d.Range($totalsParagraph.Range.End - 1, $totalsParagraph.Range.End)
$pilcrow.Delete()

# 3) Delete the now-orphaned hyperlink paragraph ("task14.xls") completely,
#    including its own paragraph mark, so the following paragraph
#    ("Выполните задание") directly follows the merged paragraph above.
$linkRange = $d.Content
$linkRange.Find.Execute("task14.xls") | Out-Null
$linkParagraph = $linkRange.Paragraphs(1)
$linkParagraph.Range.Delete()

# Sanity check: neither the old label nor the hyperlink text should remain,
# and the totals sentence should now end the (merged) paragraph.
$stillThere = $d.Content.Find.Execute("Ссылка на таблицу:")
$linkStillThere = $d.Content.Find.Execute("task14.xls")
Write-Output "leftover label found=$stillThere leftover link found=$linkStillThere"

Write-Output "done"
